$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column G
$ws.Range("G1").Value = "Calificación"
$ws.Range("G1").Font.Bold = $true

# Add "Calificación" values for rows 2-9
$ws.Range("G2").Value = 10
$ws.Range("G3").Value = 8
$ws.Range("G4").Value = 8.6
$ws.Range("G5").Value = 9.1
$ws.Range("G6").Value = 10
$ws.Range("G7").Value = 5.6
$ws.Range("G8").Value = 7.3
$ws.Range("G9").Value = 0.1

# Update selection to match the target state
$ws.Range("G12").Select()
